# NEW SAAS DATASET VERSION
# Applies the pricingData.xlsx update: refreshed figures for several
# existing rows plus three brand-new 2022 rows (Github, Jira, Overleaf).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Point fixes to existing rows (values refreshed in the new dataset
#    pull). Dependent `numberOfFeatures` (column H) cells are driven by
#    shared SUM(C:G) formulas already in the sheet, so they recalc on
#    their own once the inputs below are written.
# ---------------------------------------------------------------------

$ws.Range("I11").Value = 0

$ws.Range("G21").Value = 1
$ws.Range("K21").Value = 3

$ws.Range("K24").Value = 4

$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 23
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 7
$ws.Range("K26").Value = 5

$ws.Range("E28").Value = 40
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3

$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 33
$ws.Range("G31").Value = 1
$ws.Range("J31").Value = 22

$ws.Range("C34").Value = 1
$ws.Range("E34").Value = 6
$ws.Range("K34").Value = 3

$ws.Range("C36").Value = 2
$ws.Range("D36").Value = 3
$ws.Range("E36").Value = 8

$ws.Range("I38").Value = 11

$ws.Range("C47").Value = 1
$ws.Range("E47").Value = 7
$ws.Range("K47").Value = 4

# ---------------------------------------------------------------------
# 2. Append three brand-new rows (54-56) for 2022: Github, Jira,
#    Overleaf. Copy the formatting from the last existing data row (53)
#    first so the new rows pick up the same table styling, then fill in
#    the values/formula.
# ---------------------------------------------------------------------

$ws.Range("A53:K53").Copy()
$ws.Range("A54:K54").PasteSpecial(-4122)
$ws.Range("A55:K55").PasteSpecial(-4122)
$ws.Range("A56:K56").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("A54:K56").RowHeight = 19

$ws.Range("A54").Value = "Github"
$ws.Range("B54").Value = 2022
$ws.Range("C54").Value = 3
$ws.Range("D54").Value = 5
$ws.Range("E54").Value = 21
$ws.Range("F54").Value = 9
$ws.Range("G54").Value = 9
$ws.Range("H54").Formula = "=SUM(C54:G54)"
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 17
$ws.Range("K54").Value = 3

$ws.Range("A55").Value = "Jira"
$ws.Range("B55").Value = 2022
$ws.Range("C55").Value = 3
$ws.Range("D55").Value = 1
$ws.Range("E55").Value = 19
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 7
$ws.Range("H55").Formula = "=SUM(C55:G55)"
$ws.Range("I55").Value = 1
$ws.Range("J55").Value = 12
$ws.Range("K55").Value = 3

$ws.Range("A56").Value = "Overleaf"
$ws.Range("B56").Value = 2022
$ws.Range("C56").Value = 1
$ws.Range("D56").Value = 5
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Formula = "=SUM(C56:G56)"
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 3
$ws.Range("K56").Value = 4

# ---------------------------------------------------------------------
# 3. Grow the table / autofilter / dimension to the new extent.
# ---------------------------------------------------------------------

$tbl = $ws.ListObjects.Item("pricingData")
$tbl.Resize($ws.Range("A1:K56"))

# ---------------------------------------------------------------------
# 4. View-state tweaks captured in the sheet XML: the frozen pane now
#    starts at column B (previously G) and the lingering I11 selection
#    is cleared back to the top-left of the scrolling pane.
# ---------------------------------------------------------------------

$ws.Range("B1").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
